# Changed trade entry to limit order.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Test # 1): update From/To dates, TP %, SL %, and Strategy
$ws.Range("D2").Value = 44470
$ws.Range("E2").Value = 44561
$ws.Range("H2").Value = 9
$ws.Range("I2").Value = 6
$ws.Range("L2").Value = "EarlyMACD"

# Row 3 (Test # 2): update Exchange, From/To dates, TP %, SL %
$ws.Range("B3").Value = "Binance"
$ws.Range("D3").Value = 44470
$ws.Range("E3").Value = 44561
$ws.Range("H3").Value = 9
$ws.Range("I3").Value = 6

# Update the active cell selection to match the saved view state
$ws.Range("J7").Select()
